$wb = $excel.ActiveWorkbook

# Fix sheet name error: rename "NewSheet" to "Sheet1"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"
